$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 189 (Orégano, Mercado
# Mayorista Lo Valledor de Santiago). Insert a blank row there so every
# existing record from row 189 down shifts one row lower (189->190, ...,
# 253->254), then fill the new row with the new record's data.
$ws.Rows.Item(189).Insert()

$ws.Cells.Item(189, 1).Value  = 6
$ws.Cells.Item(189, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(189, 3).Value  = "Metropolitana"
$ws.Cells.Item(189, 4).Value  = 44876
$ws.Cells.Item(189, 5).Value  = 13
$ws.Cells.Item(189, 6).Value  = 100112029
$ws.Cells.Item(189, 7).Value  = "Orégano"
$ws.Cells.Item(189, 8).Value  = "Sin especificar"
$ws.Cells.Item(189, 9).Value  = "Primera"
$ws.Cells.Item(189, 10).Value = 48
$ws.Cells.Item(189, 11).Value = 16000
$ws.Cells.Item(189, 12).Value = 17000
$ws.Cells.Item(189, 13).Value = 16479
$ws.Cells.Item(189, 14).Value = "`$/docena de atados"
$ws.Cells.Item(189, 15).Value = "Región Metropolitana"
$ws.Cells.Item(189, 16).Value = 5493
$ws.Cells.Item(189, 17).Value = 3
$ws.Cells.Item(189, 18).Value = "Hortaliza"
